$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.796.45'
$ws.Range("E2").Value = '  +0.59%  '

$ws.Range("D3").Value = '2.622.82'
$ws.Range("E3").Value = '  -0.70%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").Value = '595.51'
$ws.Range("E5").Value = '  -0.98%  '

$ws.Range("D6").Value = '149.52'
$ws.Range("E6").Value = '  +2.23%  '

$ws.Range("E7").Value = '  -0.02%  '

$ws.Range("D8").Value = '0.587'
$ws.Range("E8").Value = '  +0.12%  '

$ws.Range("D9").Value = '0.109'
$ws.Range("E9").Value = '  +0.61%  '

$ws.Range("D10").Value = '0.382'
$ws.Range("E10").Value = '  +4.18%  '

$ws.Range("D11").Value = '5.59'
$ws.Range("E11").Value = '  -0.23%  '

$ws.Range("E12").Value = '  -1.12%  '

$ws.Range("D13").Value = '27.54'
$ws.Range("E13").Value = '  +0.84%  '

$ws.Range("D14").Value = '3.093.13'
$ws.Range("E14").Value = '  -1.01%  '

$ws.Range("D15").Value = '63.661.94'
$ws.Range("E15").Value = '  +0.58%  '

$ws.Range("E16").Value = '  +2.24%  '

$ws.Range("D17").Value = '2.621.45'
$ws.Range("E17").Value = '  -0.32%  '

$ws.Range("D18").Value = '12.18'
$ws.Range("E18").Value = '  +6.47%  '

$ws.Range("D19").Value = '4.61'
$ws.Range("E19").Value = '  +1.99%  '

$ws.Range("D20").Value = '348.58'
$ws.Range("E20").Value = '  +2.16%  '

$ws.Range("D21").Value = '6.86'
$ws.Range("E21").Value = '  -0.52%  '

$ws.Range("D22").Value = '0.999'
$ws.Range("E22").Value = '  -0.07%  '

$ws.Range("E23").Value = '  +2.52%  '

$ws.Range("D24").Value = '66.11'
$ws.Range("E24").Value = '  -0.98%  '

$ws.Range("D25").Value = '1.73'
$ws.Range("E25").Value = '  +13.00%  '

$ws.Range("E26").Value = '  -0.68%  '

$ws.Range("D27").Value = '9.20'
$ws.Range("E27").Value = '  +3.80%  '

$ws.Range("D28").Value = '0.164'
$ws.Range("E28").Value = '  +0.28%  '

$ws.Range("D29").Value = '8.08'
$ws.Range("E29").Value = '  +2.88%  '

$ws.Range("D30").Value = '543.73'
$ws.Range("E30").Value = '  -0.48%  '

$ws.Range("E31").Value = '  +0.07%  '

$ws.Range("D32").Value = '2.03'
$ws.Range("E32").Value = '  -1.14%  '

$ws.Range("D33").Value = '0.0₃0844'
$ws.Range("E33").Value = '  +4.90%  '

$ws.Range("D34").Value = '1.75'
$ws.Range("E34").Value = '  +0.13%  '

$ws.Range("D35").Value = '5.21'
$ws.Range("E35").Value = '  -0.32%  '

$ws.Range("D36").Value = '168.01'
$ws.Range("E36").Value = '  +0.03%  '

$ws.Range("D37").Value = '0.406'
$ws.Range("E37").Value = '  +0.23%  '

$ws.Range("D38").Value = '1.00'
$ws.Range("E38").Value = '  -0.07%  '

$ws.Range("D39").Value = '1.95'
$ws.Range("E39").Value = '  +2.72%  '

$ws.Range("D40").Value = '19.38'
$ws.Range("E40").Value = '  +1.61%  '

$ws.Range("E41").Value = '  +0.13%  '

$ws.Range("D42").Value = '169.73'
$ws.Range("E42").Value = '  +0.44%  '

$ws.Range("D43").Value = '39.79'
$ws.Range("E43").Value = '  +0.13%  '

$ws.Range("D44").Value = '3.93'
$ws.Range("E44").Value = '  +4.16%  '

$ws.Range("D45").Value = '0.0593'
$ws.Range("E45").Value = '  +2.91%  '

$ws.Range("D46").Value = '21.39'
$ws.Range("E46").Value = '  -4.91%  '

$ws.Range("D47").Value = '0.627'
$ws.Range("E47").Value = '  +0.35%  '

$ws.Range("B48").Value = 'VeChain'
$ws.Range("C48").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D48").Value = '0.0245'
$ws.Range("E48").Value = '  -0.65%  '

$ws.Range("B49").Value = 'dogwifhat'
$ws.Range("C49").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D49").Value = '1.98'
$ws.Range("E49").Value = '  +10.70%  '

$ws.Range("D50").Value = '0.0966'
$ws.Range("E50").Value = '  +0.37%  '

$ws.Range("D51").Value = '19.19'
$ws.Range("E51").Value = '  +1.93%  '

